# Apply the "Lloyd shipping economist" update:
#  - add a new shared-string value "newbuilding_price_per_dwt_2500teu_fullcon"
#  - change column D (type) for rows 10-20 from the 1600teu string to the new 2500teu string
#  - fill in missing E (price_main) values for rows 13-20 and 32-39
#  - adjust sheet view (scroll position / selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update the "type" column (D) for rows 10 through 20 to the new category ---
$newType = "newbuilding_price_per_dwt_2500teu_fullcon"
for ($r = 10; $r -le 20; $r++) {
    $ws.Range("D$r").Value = $newType
}

# --- Fill in the previously-empty price_main values (column E) ---
$values = @{
    13 = 1616
    14 = 1282
    15 = 1205
    16 = 1077
    17 = 2000
    18 = 1800
    19 = 1800
    20 = 1150
    32 = 865
    33 = 816
    34 = 790
    35 = 827
    36 = 977
    37 = 977
    38 = 977
    39 = 751
}

foreach ($row in $values.Keys) {
    $ws.Range("E$row").Value = $values[$row]
}

# --- Adjust the sheet view: scroll back to top-left and move the selection ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E40").Select()
